# Refresh the crypto price/volume snapshot (and the two ranking swaps
# at rows 32/33 and 43/44) to match the latest coinranking.com pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 carries the sheet's plain ("Normal") style. We use it to wipe any
# incidental number-format styling that Excel applies when we briefly
# force a cell to Text so numeric-looking values (e.g. "1.001") are not
# silently reinterpreted as numbers - the source data is text throughout.
$normalStyle = $ws.Range("B2").Style

$updates = @(
    @{ Addr = 'D2'; Value = '23.481.45' }
    @{ Addr = 'E2'; Value = '  +0.34%  ' }
    @{ Addr = 'D3'; Value = '1.640.40' }
    @{ Addr = 'E3'; Value = '  -0.20%  ' }
    @{ Addr = 'D4'; Value = '1.001' }
    @{ Addr = 'E4'; Value = '  -0.85%  ' }
    @{ Addr = 'D5'; Value = '0.9987' }
    @{ Addr = 'E5'; Value = '  -0.53%  ' }
    @{ Addr = 'D6'; Value = '304.23' }
    @{ Addr = 'E6'; Value = '  +0.04%  ' }
    @{ Addr = 'D7'; Value = '0.3792' }
    @{ Addr = 'E7'; Value = '  +0.79%  ' }
    @{ Addr = 'D8'; Value = '51.69' }
    @{ Addr = 'E8'; Value = '  -1.92%  ' }
    @{ Addr = 'E9'; Value = '  -0.71%  ' }
    @{ Addr = 'D10'; Value = '0.08179' }
    @{ Addr = 'E10'; Value = '  +0.67%  ' }
    @{ Addr = 'E11'; Value = '  -1.43%  ' }
    @{ Addr = 'D12'; Value = '0.9983' }
    @{ Addr = 'E12'; Value = '  -1.09%  ' }
    @{ Addr = 'D13'; Value = '22.55' }
    @{ Addr = 'E13'; Value = '  -1.65%  ' }
    @{ Addr = 'D14'; Value = '6.476' }
    @{ Addr = 'E14'; Value = '  -2.75%  ' }
    @{ Addr = 'D15'; Value = '7.383' }
    @{ Addr = 'E15'; Value = '  +0.61%  ' }
    @{ Addr = 'D16'; Value = '0.00001241' }
    @{ Addr = 'E16'; Value = '  -1.69%  ' }
    @{ Addr = 'D17'; Value = '1.637.28' }
    @{ Addr = 'E17'; Value = '  -0.69%  ' }
    @{ Addr = 'D18'; Value = '95.62' }
    @{ Addr = 'E18'; Value = '  +1.21%  ' }
    @{ Addr = 'D19'; Value = '0.06929' }
    @{ Addr = 'E19'; Value = '  +0.10%  ' }
    @{ Addr = 'D20'; Value = '6.596' }
    @{ Addr = 'E20'; Value = '  +0.32%  ' }
    @{ Addr = 'D21'; Value = '17.54' }
    @{ Addr = 'E21'; Value = '  -3.87%  ' }
    @{ Addr = 'D22'; Value = '0.9983' }
    @{ Addr = 'E22'; Value = '  -0.42%  ' }
    @{ Addr = 'D23'; Value = '12.52' }
    @{ Addr = 'E23'; Value = '  -2.93%  ' }
    @{ Addr = 'D24'; Value = '23.490.44' }
    @{ Addr = 'E24'; Value = '  +0.25%  ' }
    @{ Addr = 'D25'; Value = '2.505' }
    @{ Addr = 'E25'; Value = '  +2.98%  ' }
    @{ Addr = 'D26'; Value = '3.064' }
    @{ Addr = 'E26'; Value = '  -4.98%  ' }
    @{ Addr = 'D27'; Value = '21.18' }
    @{ Addr = 'E27'; Value = '  -0.76%  ' }
    @{ Addr = 'D28'; Value = '151.63' }
    @{ Addr = 'E28'; Value = '  -0.12%  ' }
    @{ Addr = 'D29'; Value = '5.238' }
    @{ Addr = 'E29'; Value = '  -0.99%  ' }
    @{ Addr = 'D30'; Value = '133.40' }
    @{ Addr = 'E30'; Value = '  -2.52%  ' }
    @{ Addr = 'D31'; Value = '1.818.79' }
    @{ Addr = 'E31'; Value = '  -0.70%  ' }
    @{ Addr = 'B32'; Value = 'Filecoin' }
    @{ Addr = 'C32'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Addr = 'D32'; Value = '6.652' }
    @{ Addr = 'E32'; Value = '  -4.27%  ' }
    @{ Addr = 'B33'; Value = 'WEMIXTOKEN' }
    @{ Addr = 'C33'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' }
    @{ Addr = 'D33'; Value = '2.158' }
    @{ Addr = 'E33'; Value = '  -7.01%  ' }
    @{ Addr = 'D34'; Value = '1.072' }
    @{ Addr = 'E34'; Value = '  +11.18%  ' }
    @{ Addr = 'D35'; Value = '11.47' }
    @{ Addr = 'E35'; Value = '  +4.83%  ' }
    @{ Addr = 'D36'; Value = '0.02766' }
    @{ Addr = 'E36'; Value = '  -3.43%  ' }
    @{ Addr = 'D37'; Value = '0.2496' }
    @{ Addr = 'E37'; Value = '  -2.85%  ' }
    @{ Addr = 'D38'; Value = '0.08778' }
    @{ Addr = 'E38'; Value = '  -0.93%  ' }
    @{ Addr = 'D39'; Value = '0.07111' }
    @{ Addr = 'E39'; Value = '  -2.48%  ' }
    @{ Addr = 'D40'; Value = '6.026' }
    @{ Addr = 'E40'; Value = '  -4.77%  ' }
    @{ Addr = 'D41'; Value = '0.7058' }
    @{ Addr = 'E41'; Value = '  -1.43%  ' }
    @{ Addr = 'D42'; Value = '1.347' }
    @{ Addr = 'E42'; Value = '  -2.42%  ' }
    @{ Addr = 'B43'; Value = 'Aptos' }
    @{ Addr = 'C43'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Addr = 'D43'; Value = '12.22' }
    @{ Addr = 'E43'; Value = '  -3.44%  ' }
    @{ Addr = 'B44'; Value = 'EnergySwap' }
    @{ Addr = 'C44'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Addr = 'D44'; Value = '15.73' }
    @{ Addr = 'E44'; Value = '  -5.38%  ' }
    @{ Addr = 'D45'; Value = '0.6542' }
    @{ Addr = 'E45'; Value = '  -0.82%  ' }
    @{ Addr = 'D46'; Value = '0.9982' }
    @{ Addr = 'E46'; Value = '  -0.34%  ' }
    @{ Addr = 'D47'; Value = '2.286' }
    @{ Addr = 'E47'; Value = '  -3.34%  ' }
    @{ Addr = 'D48'; Value = '3.972' }
    @{ Addr = 'E48'; Value = '  -1.00%  ' }
    @{ Addr = 'D49'; Value = '0.07981' }
    @{ Addr = 'E49'; Value = '  -0.38%  ' }
    @{ Addr = 'D50'; Value = '128.76' }
    @{ Addr = 'E50'; Value = '  +0.57%  ' }
    @{ Addr = 'D51'; Value = '1.194' }
    @{ Addr = 'E51'; Value = '  -2.09%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = $normalStyle
}

